$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf2"
$ws.Range("C2").Value = "Cd44"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2151743333333333
$ws.Range("H2").Value = 0.645523
$ws.Range("I2").Value = 0.01945888736810517
$ws.Range("J2").Value = 0.01945888736810517
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 239.0839323333333
$ws.Range("N2").Value = 717.251797
$ws.Range("O2").Value = 0.4086975387666237
$ws.Range("P2").Value = 0.4086975387666237
$ws.Range("Q2").Value = 51.44472575053678
$ws.Range("R2").Value = 463.002531754831
$ws.Range("S2").Value = 0.007952799374481526
$ws.Range("T2").Value = 0.007952799374481526

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf2"
$ws.Range("C3").Value = "Cd44"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2151743333333333
$ws.Range("H3").Value = 0.645523
$ws.Range("I3").Value = 0.01945888736810517
$ws.Range("J3").Value = 0.01945888736810517
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 117.0512696666667
$ws.Range("N3").Value = 351.153809
$ws.Range("O3").Value = 0.2000910950200451
$ws.Range("P3").Value = 0.2000910950200451
$ws.Range("Q3").Value = 25.18642891634522
$ws.Range("R3").Value = 226.677860247107
$ws.Range("S3").Value = 0.003893550081355886
$ws.Range("T3").Value = 0.003893550081355886

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf2"
$ws.Range("C4").Value = "Cd44"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2151743333333333
$ws.Range("H4").Value = 0.645523
$ws.Range("I4").Value = 0.01945888736810517
$ws.Range("J4").Value = 0.01945888736810517
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 171.15883
$ws.Range("N4").Value = 513.47649
$ws.Range("O4").Value = 0.2925842480357353
$ws.Range("P4").Value = 0.2925842480357353
$ws.Range("Q4").Value = 36.82898713936333
$ws.Range("R4").Value = 331.46088425427
$ws.Range("S4").Value = 0.005693363928209119
$ws.Range("T4").Value = 0.005693363928209119

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf2"
$ws.Range("C5").Value = "Cd44"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.2151743333333333
$ws.Range("H5").Value = 0.645523
$ws.Range("I5").Value = 0.01945888736810517
$ws.Range("J5").Value = 0.01945888736810517
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 57.695868
$ws.Range("N5").Value = 173.087604
$ws.Range("O5").Value = 0.09862711817759588
$ws.Range("P5").Value = 0.09862711817759588
$ws.Range("Q5").Value = 12.414669932988
$ws.Range("R5").Value = 111.732029396892
$ws.Range("S5").Value = 0.001919173984058636
$ws.Range("T5").Value = 0.001919173984058636

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fgf2"
$ws.Range("C6").Value = "Cd44"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.077181333333334
$ws.Range("H6").Value = 24.231544
$ws.Range("I6").Value = 0.730444748601188
$ws.Range("J6").Value = 0.730444748601188
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 239.0839323333333
$ws.Range("N6").Value = 717.251797
$ws.Range("O6").Value = 0.4086975387666237
$ws.Range("P6").Value = 0.4086975387666237
$ws.Range("Q6").Value = 1931.12427534273
$ws.Range("R6").Value = 17380.11847808457
$ws.Range("S6").Value = 0.2985309709583107
$ws.Range("T6").Value = 0.2985309709583107

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fgf2"
$ws.Range("C7").Value = "Cd44"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 8.077181333333334
$ws.Range("H7").Value = 24.231544
$ws.Range("I7").Value = 0.730444748601188
$ws.Range("J7").Value = 0.730444748601188
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 117.0512696666667
$ws.Range("N7").Value = 351.153809
$ws.Range("O7").Value = 0.2000910950200451
$ws.Range("P7").Value = 0.2000910950200451
$ws.Range("Q7").Value = 945.4443303945664
$ws.Range("R7").Value = 8508.998973551097
$ws.Range("S7").Value = 0.1461554895992533
$ws.Range("T7").Value = 0.1461554895992533

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fgf2"
$ws.Range("C8").Value = "Cd44"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 8.077181333333334
$ws.Range("H8").Value = 24.231544
$ws.Range("I8").Value = 0.730444748601188
$ws.Range("J8").Value = 0.730444748601188
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 171.15883
$ws.Range("N8").Value = 513.47649
$ws.Range("O8").Value = 0.2925842480357353
$ws.Range("P8").Value = 0.2925842480357353
$ws.Range("Q8").Value = 1382.480906711173
$ws.Range("R8").Value = 12442.32816040056
$ws.Range("S8").Value = 0.2137166275011303
$ws.Range("T8").Value = 0.2137166275011303

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fgf2"
$ws.Range("C9").Value = "Cd44"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 8.077181333333334
$ws.Range("H9").Value = 24.231544
$ws.Range("I9").Value = 0.730444748601188
$ws.Range("J9").Value = 0.730444748601188
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 57.695868
$ws.Range("N9").Value = 173.087604
$ws.Range("O9").Value = 0.09862711817759588
$ws.Range("P9").Value = 0.09862711817759588
$ws.Range("Q9").Value = 466.019988020064
$ws.Range("R9").Value = 4194.179892180576
$ws.Range("S9").Value = 0.07204166054249368
$ws.Range("T9").Value = 0.07204166054249368

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Fgf2"
$ws.Range("C10").Value = "Cd44"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.519176
$ws.Range("H10").Value = 1.557528
$ws.Range("I10").Value = 0.04695070806875992
$ws.Range("J10").Value = 0.04695070806875992
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 239.0839323333333
$ws.Range("N10").Value = 717.251797
$ws.Range("O10").Value = 0.4086975387666237
$ws.Range("P10").Value = 0.4086975387666237
$ws.Range("Q10").Value = 124.1266396530907
$ws.Range("R10").Value = 1117.139756877816
$ws.Range("S10").Value = 0.01918863883105244
$ws.Range("T10").Value = 0.01918863883105244

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Fgf2"
$ws.Range("C11").Value = "Cd44"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.519176
$ws.Range("H11").Value = 1.557528
$ws.Range("I11").Value = 0.04695070806875992
$ws.Range("J11").Value = 0.04695070806875992
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 117.0512696666667
$ws.Range("N11").Value = 351.153809
$ws.Range("O11").Value = 0.2000910950200451
$ws.Range("P11").Value = 0.2000910950200451
$ws.Range("Q11").Value = 60.77020998046133
$ws.Range("R11").Value = 546.931889824152
$ws.Range("S11").Value = 0.009394418589444637
$ws.Range("T11").Value = 0.009394418589444639

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Fgf2"
$ws.Range("C12").Value = "Cd44"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.519176
$ws.Range("H12").Value = 1.557528
$ws.Range("I12").Value = 0.04695070806875992
$ws.Range("J12").Value = 0.04695070806875992
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 171.15883
$ws.Range("N12").Value = 513.47649
$ws.Range("O12").Value = 0.2925842480357353
$ws.Range("P12").Value = 0.2925842480357353
$ws.Range("Q12").Value = 88.86155672407999
$ws.Range("R12").Value = 799.7540105167201
$ws.Range("S12").Value = 0.01373703761504345
$ws.Range("T12").Value = 0.01373703761504345

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Fgf2"
$ws.Range("C13").Value = "Cd44"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.519176
$ws.Range("H13").Value = 1.557528
$ws.Range("I13").Value = 0.04695070806875992
$ws.Range("J13").Value = 0.04695070806875992
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 57.695868
$ws.Range("N13").Value = 173.087604
$ws.Range("O13").Value = 0.09862711817759588
$ws.Range("P13").Value = 0.09862711817759588
$ws.Range("Q13").Value = 29.954309964768
$ws.Range("R13").Value = 269.588789682912
$ws.Range("S13").Value = 0.004630613033219389
$ws.Range("T13").Value = 0.004630613033219389

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Fgf2"
$ws.Range("C14").Value = "Cd44"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.246363333333333
$ws.Range("H14").Value = 6.73909
$ws.Range("I14").Value = 0.2031456559619469
$ws.Range("J14").Value = 0.2031456559619469
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 239.0839323333333
$ws.Range("N14").Value = 717.251797
$ws.Range("O14").Value = 0.4086975387666237
$ws.Range("P14").Value = 0.4086975387666237
$ws.Range("Q14").Value = 537.0693791827479
$ws.Range("R14").Value = 4833.62441264473
$ws.Range("S14").Value = 0.08302512960277901
$ws.Range("T14").Value = 0.08302512960277901

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Fgf2"
$ws.Range("C15").Value = "Cd44"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.246363333333333
$ws.Range("H15").Value = 6.73909
$ws.Range("I15").Value = 0.2031456559619469
$ws.Range("J15").Value = 0.2031456559619469
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 117.0512696666667
$ws.Range("N15").Value = 351.153809
$ws.Range("O15").Value = 0.2000910950200451
$ws.Range("P15").Value = 0.2000910950200451
$ws.Range("Q15").Value = 262.9396802993123
$ws.Range("R15").Value = 2366.45712269381
$ws.Range("S15").Value = 0.04064763674999131
$ws.Range("T15").Value = 0.04064763674999131

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Fgf2"
$ws.Range("C16").Value = "Cd44"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.246363333333333
$ws.Range("H16").Value = 6.73909
$ws.Range("I16").Value = 0.2031456559619469
$ws.Range("J16").Value = 0.2031456559619469
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 171.15883
$ws.Range("N16").Value = 513.47649
$ws.Range("O16").Value = 0.2925842480357353
$ws.Range("P16").Value = 0.2925842480357353
$ws.Range("Q16").Value = 384.4849198882333
$ws.Range("R16").Value = 3460.3642789941
$ws.Range("S16").Value = 0.05943721899135243
$ws.Range("T16").Value = 0.05943721899135243

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Fgf2"
$ws.Range("C17").Value = "Cd44"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2.246363333333333
$ws.Range("H17").Value = 6.73909
$ws.Range("I17").Value = 0.2031456559619469
$ws.Range("J17").Value = 0.2031456559619469
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 57.695868
$ws.Range("N17").Value = 173.087604
$ws.Range("O17").Value = 0.09862711817759588
$ws.Range("P17").Value = 0.09862711817759588
$ws.Range("Q17").Value = 129.60588236004
$ws.Range("R17").Value = 1166.45294124036
$ws.Range("S17").Value = 0.02003567061782418
$ws.Range("T17").Value = 0.02003567061782418
